$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing B-column values (rows 2-127) for the rerun on the new dataset ---
$bUpdates = @{
    2 = 0.28864
    3 = 0.28864
    4 = 0.27809
    5 = 0.27809
    6 = 0.27809
    7 = 0.27809
    8 = 0.27685
    9 = 0.27809
    10 = 0.27498
    11 = 0.2905
    12 = 0.28305
    13 = 0.28305
    14 = 0.28305
    15 = 0.28305
    16 = 0.27312
    17 = 0.27498
    18 = 0.27685
    19 = 0.26567
    20 = 0.2905
    21 = 0.2725
    22 = 0.28802
    23 = 0.2725
    24 = 0.2843
    25 = 0.28367
    26 = 0.2843
    27 = 0.27312
    28 = 0.27809
    29 = 0.2843
    30 = 0.27809
    31 = 0.26754
    32 = 0.27747
    33 = 0.28181
    35 = 0.28988
    36 = 0.28988
    37 = 0.28988
    38 = 0.29981
    39 = 0.29671
    40 = 0.28988
    41 = 0.28926
    42 = 0.27747
    43 = 0.27312
    44 = 0.27312
    45 = 0.28181
    46 = 0.27623
    47 = 0.27747
    48 = 0.27623
    49 = 0.27623
    50 = 0.26567
    51 = 0.27064
    52 = 0.26878
    53 = 0.26567
    54 = 0.25947
    55 = 0.25885
    56 = 0.26567
    57 = 0.25947
    58 = 0.2545
    59 = 0.25885
    60 = 0.27498
    61 = 0.27126
    62 = 0.26505
    63 = 0.25947
    64 = 0.26505
    65 = 0.26071
    66 = 0.26133
    67 = 0.26133
    68 = 0.27188
    69 = 0.26816
    70 = 0.28367
    71 = 0.2874
    72 = 0.28988
    73 = 0.26691
    74 = 0.28119
    75 = 0.27126
    76 = 0.26878
    77 = 0.28181
    78 = 0.28181
    79 = 0.28367
    80 = 0.30292
    81 = 0.29361
    82 = 0.29671
    83 = 0.30912
    84 = 0.31533
    85 = 0.31844
    86 = 0.30912
    87 = 0.29671
    88 = 0.29671
    90 = 0.2905
    91 = 0.30292
    92 = 0.29671
    93 = 0.29671
    94 = 0.2905
    95 = 0.29361
    96 = 0.28492
    97 = 0.28243
    98 = 0.27188
    99 = 0.26816
    100 = 0.2843
    101 = 0.30602
    102 = 0.30602
    103 = 0.29671
    104 = 0.28802
    105 = 0.2843
    106 = 0.27933
    107 = 0.28926
    108 = 0.2905
    109 = 0.30292
    110 = 0.30912
    111 = 0.29671
    112 = 0.28554
    113 = 0.2905
    114 = 0.31533
    115 = 0.30292
    116 = 0.30912
    117 = 0.29361
    118 = 0.31533
    119 = 0.35878
    120 = 0.39603
    121 = 0.39603
    122 = 0.38361
    123 = 0.38982
    124 = 0.38051
    125 = 0.38051
    126 = 0.37741
    127 = 0.40534
}

foreach ($row in $bUpdates.Keys) {
    $ws.Cells.Item($row, 2).Value = $bUpdates[$row]
}

# --- Append new rows 128-133 (A: index, B: close value) ---
$newRows = @(
    @{ Row = 128; A = 126; B = 0.3712 }
    @{ Row = 129; A = 127; B = 0.3712 }
    @{ Row = 130; A = 128; B = 0.41775 }
    @{ Row = 131; A = 129; B = 0.45189 }
    @{ Row = 132; A = 130; B = 0.44569 }
    @{ Row = 133; A = 131; B = 0.4581 }
)

# Copy the format from the last existing row (127) down onto the new rows so
# column A keeps its centered/bordered style, then set the actual values.
$ws.Range("A127").Copy() | Out-Null
$ws.Range("A128:A133").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
}

